$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the previously-missing values in row 3 (H3 and I3)
$ws.Range("H3").Value = 378
$ws.Range("I3").Value = 54

# Update the view: scroll back to A1 (remove topLeftCell offset) and move
# the active selection to N9
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("N9").Select()
